$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that would otherwise be auto-parsed as numbers by Excel need to be
# explicitly formatted as Text first so the literal string is preserved.
$textCells = @("D5", "D8", "D13", "D16", "D18", "D19", "D21", "D23", "D24", "D26", "D30", "D32", "D38", "D40", "D41", "D44", "D47", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.073.40"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "1.567.84"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("E4").Value = "  +0.58%  "
$ws.Range("D5").Value = "208.75"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("E7").Value = "  +0.65%  "
$ws.Range("D8").Value = "22.04"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("E10").Value = "  +1.44%  "
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").Value = "1.564.15"
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("D13").Value = "3.77"
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("D15").Value = "27.065.10"
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").Value = "61.96"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D18").Value = "7.41"
$ws.Range("E18").Value = "  +2.06%  "
$ws.Range("D19").Value = "215.45"
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("D21").Value = "4.14"
$ws.Range("E21").Value = "  +2.03%  "
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("D23").Value = "1.95"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "154.05"
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("D26").Value = "15.04"
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("E27").Value = "  +1.28%  "
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("E29").Value = "  +4.45%  "
$ws.Range("D30").Value = "0.0472"
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("D32").Value = "3.20"
$ws.Range("E32").Value = "  +2.99%  "
$ws.Range("D33").Value = "1.428.42"
$ws.Range("E33").Value = "  +1.30%  "
$ws.Range("E34").Value = "  +13.52%  "
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("E36").Value = "  +3.42%  "
$ws.Range("E37").Value = "  +1.19%  "
$ws.Range("D38").Value = "0.531"
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("E39").Value = "  +2.71%  "
$ws.Range("D40").Value = "0.812"
$ws.Range("E40").Value = "  +0.75%  "
$ws.Range("D41").Value = "2.37"
$ws.Range("E41").Value = "  +3.71%  "
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("E43").Value = "  +0.58%  "
$ws.Range("D44").Value = "64.65"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").Value = "1.704.41"
$ws.Range("E46").Value = "  +1.16%  "
$ws.Range("D47").Value = "86.72"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("E48").Value = "  +3.17%  "
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").Value = "0.0963"
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("E51").Value = "  +0.64%  "
